# Auto: Update ETF Data
# Applies the diff: adds Close Price / Market Value / Share Change / Net Amount
# columns (E:H) to the holdings sheet, renames A1 header, and extends the
# used dimension range accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Data Date" label in A1 to "Date"
$ws.Range("A1").Value = "Date"

# Copy the existing header style (from D3, "Weight") onto the new header
# cells E3:H3 before writing their text, so they pick up the bold / bordered
# / centered formatting used for the other column headers.
$ws.Range("D3").Copy()
$ws.Range("E3:H3").PasteSpecial(-4122)

$ws.Range("E3").Value = "Close Price"
$ws.Range("F3").Value = "Market Value"
$ws.Range("G3").Value = "Share Change"
$ws.Range("H3").Value = "Net Amount"

$ws.Cells.Item(4, 5).Value = 1550
$ws.Cells.Item(4, 6).Value = 4876300000
$ws.Cells.Item(4, 7).Value = 0
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(5, 5).Value = 1645
$ws.Cells.Item(5, 6).Value = 3350865000
$ws.Cells.Item(5, 7).Value = 0
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(6, 5).Value = 1510
$ws.Cells.Item(6, 6).Value = 3225360000
$ws.Cells.Item(6, 7).Value = 0
$ws.Cells.Item(6, 8).Value = 0
$ws.Cells.Item(7, 5).Value = 4485
$ws.Cells.Item(7, 6).Value = 3108105000
$ws.Cells.Item(7, 7).Value = 0
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(8, 5).Value = 687
$ws.Cells.Item(8, 6).Value = 2910819000
$ws.Cells.Item(8, 7).Value = 0
$ws.Cells.Item(8, 8).Value = 0
$ws.Cells.Item(9, 5).Value = 1185
$ws.Cells.Item(9, 6).Value = 2793045000
$ws.Cells.Item(9, 7).Value = 0
$ws.Cells.Item(9, 8).Value = 0
$ws.Cells.Item(10, 5).Value = 963
$ws.Cells.Item(10, 6).Value = 2438316000
$ws.Cells.Item(10, 7).Value = 0
$ws.Cells.Item(10, 8).Value = 0
$ws.Cells.Item(11, 5).Value = 2250
$ws.Cells.Item(11, 6).Value = 2423250000
$ws.Cells.Item(11, 7).Value = 0
$ws.Cells.Item(11, 8).Value = 0
$ws.Cells.Item(12, 5).Value = 494
$ws.Cells.Item(12, 6).Value = 2384538000
$ws.Cells.Item(12, 7).Value = 0
$ws.Cells.Item(12, 8).Value = 0
$ws.Cells.Item(13, 5).Value = 1520
$ws.Cells.Item(13, 6).Value = 2329928960
$ws.Cells.Item(13, 7).Value = 0
$ws.Cells.Item(13, 8).Value = 0
$ws.Cells.Item(14, 5).Value = 2745
$ws.Cells.Item(14, 6).Value = 2039535000
$ws.Cells.Item(14, 7).Value = 0
$ws.Cells.Item(14, 8).Value = 0
$ws.Cells.Item(15, 5).Value = 230.5
$ws.Cells.Item(15, 6).Value = 1867050000
$ws.Cells.Item(15, 7).Value = 0
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(16, 5).Value = 1430
$ws.Cells.Item(16, 6).Value = 1670240000
$ws.Cells.Item(16, 7).Value = 0
$ws.Cells.Item(16, 8).Value = 0
$ws.Cells.Item(17, 5).Value = 1625
$ws.Cells.Item(17, 6).Value = 1618500000
$ws.Cells.Item(17, 7).Value = 0
$ws.Cells.Item(17, 8).Value = 0
$ws.Cells.Item(18, 5).Value = 1005
$ws.Cells.Item(18, 6).Value = 1485390000
$ws.Cells.Item(18, 7).Value = 0
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(19, 5).Value = 3510
$ws.Cells.Item(19, 6).Value = 1432080000
$ws.Cells.Item(19, 7).Value = 0
$ws.Cells.Item(19, 8).Value = 0
$ws.Cells.Item(20, 5).Value = 1450
$ws.Cells.Item(20, 6).Value = 1399250000
$ws.Cells.Item(20, 7).Value = 0
$ws.Cells.Item(20, 8).Value = 0
$ws.Cells.Item(21, 5).Value = 247.5
$ws.Cells.Item(21, 6).Value = 995940000
$ws.Cells.Item(21, 7).Value = 0
$ws.Cells.Item(21, 8).Value = 0
$ws.Cells.Item(22, 5).Value = 7260
$ws.Cells.Item(22, 6).Value = 972840000
$ws.Cells.Item(22, 7).Value = 0
$ws.Cells.Item(22, 8).Value = 0
$ws.Cells.Item(23, 5).Value = 3750
$ws.Cells.Item(23, 6).Value = 870000000
$ws.Cells.Item(23, 7).Value = 0
$ws.Cells.Item(23, 8).Value = 0
$ws.Cells.Item(24, 5).Value = 949
$ws.Cells.Item(24, 6).Value = 745914000
$ws.Cells.Item(24, 7).Value = 0
$ws.Cells.Item(24, 8).Value = 0
$ws.Cells.Item(25, 5).Value = 277
$ws.Cells.Item(25, 6).Value = 462036000
$ws.Cells.Item(25, 7).Value = 0
$ws.Cells.Item(25, 8).Value = 0
$ws.Cells.Item(26, 5).Value = 93.09999847412109
$ws.Cells.Item(26, 6).Value = 436173492.8512573
$ws.Cells.Item(26, 7).Value = 0
$ws.Cells.Item(26, 8).Value = 0
$ws.Cells.Item(27, 5).Value = 781
$ws.Cells.Item(27, 6).Value = 399872000
$ws.Cells.Item(27, 7).Value = 0
$ws.Cells.Item(27, 8).Value = 0
$ws.Cells.Item(28, 5).Value = 112
$ws.Cells.Item(28, 6).Value = 345856000
$ws.Cells.Item(28, 7).Value = 0
$ws.Cells.Item(28, 8).Value = 0
$ws.Cells.Item(29, 5).Value = 1295
$ws.Cells.Item(29, 6).Value = 326340000
$ws.Cells.Item(29, 7).Value = 0
$ws.Cells.Item(29, 8).Value = 0
$ws.Cells.Item(30, 5).Value = 60.20000076293945
$ws.Cells.Item(30, 6).Value = 311414603.9466858
$ws.Cells.Item(30, 7).Value = 0
$ws.Cells.Item(30, 8).Value = 0
$ws.Cells.Item(31, 5).Value = 127.5
$ws.Cells.Item(31, 6).Value = 252827400
$ws.Cells.Item(31, 7).Value = 0
$ws.Cells.Item(31, 8).Value = 0
$ws.Cells.Item(32, 5).Value = 142
$ws.Cells.Item(32, 6).Value = 246796000
$ws.Cells.Item(32, 7).Value = 0
$ws.Cells.Item(32, 8).Value = 0
$ws.Cells.Item(33, 5).Value = 150.5
$ws.Cells.Item(33, 6).Value = 233275000
$ws.Cells.Item(33, 7).Value = 0
$ws.Cells.Item(33, 8).Value = 0
$ws.Cells.Item(34, 5).Value = 248
$ws.Cells.Item(34, 6).Value = 211792000
$ws.Cells.Item(34, 7).Value = 0
$ws.Cells.Item(34, 8).Value = 0
$ws.Cells.Item(35, 5).Value = 61.70000076293945
$ws.Cells.Item(35, 6).Value = 197995302.4482727
$ws.Cells.Item(35, 7).Value = 0
$ws.Cells.Item(35, 8).Value = 0
$ws.Cells.Item(36, 5).Value = 337
$ws.Cells.Item(36, 6).Value = 187035000
$ws.Cells.Item(36, 7).Value = 0
$ws.Cells.Item(36, 8).Value = 0
$ws.Cells.Item(37, 5).Value = 97
$ws.Cells.Item(37, 6).Value = 152678000
$ws.Cells.Item(37, 7).Value = 0
$ws.Cells.Item(37, 8).Value = 0
$ws.Cells.Item(38, 5).Value = 2840
$ws.Cells.Item(38, 6).Value = 136320000
$ws.Cells.Item(38, 7).Value = 0
$ws.Cells.Item(38, 8).Value = 0
$ws.Cells.Item(39, 5).Value = 92.09999847412109
$ws.Cells.Item(39, 6).Value = 88139698.53973389
$ws.Cells.Item(39, 7).Value = 0
$ws.Cells.Item(39, 8).Value = 0
$ws.Cells.Item(40, 5).Value = 159
$ws.Cells.Item(40, 6).Value = 48654000
$ws.Cells.Item(40, 7).Value = 0
$ws.Cells.Item(40, 8).Value = 0
$ws.Cells.Item(41, 5).Value = 317.5
$ws.Cells.Item(41, 6).Value = 45402500
$ws.Cells.Item(41, 7).Value = 0
$ws.Cells.Item(41, 8).Value = 0
$ws.Cells.Item(42, 5).Value = 2285
$ws.Cells.Item(42, 6).Value = 4570000
$ws.Cells.Item(42, 7).Value = 0
$ws.Cells.Item(42, 8).Value = 0
$ws.Cells.Item(43, 5).Value = 548
$ws.Cells.Item(43, 6).Value = 548000
$ws.Cells.Item(43, 7).Value = 0
$ws.Cells.Item(43, 8).Value = 0
$ws.Cells.Item(44, 5).Value = 580
$ws.Cells.Item(44, 6).Value = 580000
$ws.Cells.Item(44, 7).Value = 0
$ws.Cells.Item(44, 8).Value = 0
$ws.Cells.Item(45, 5).Value = 97.69999694824219
$ws.Cells.Item(45, 6).Value = 97699.99694824219
$ws.Cells.Item(45, 7).Value = 0
$ws.Cells.Item(45, 8).Value = 0
$ws.Cells.Item(46, 5).Value = 2495
$ws.Cells.Item(46, 6).Value = 2495000
$ws.Cells.Item(46, 7).Value = 0
$ws.Cells.Item(46, 8).Value = 0
$ws.Cells.Item(47, 5).Value = 108.5
$ws.Cells.Item(47, 6).Value = 108500
$ws.Cells.Item(47, 7).Value = 0
$ws.Cells.Item(47, 8).Value = 0
$ws.Cells.Item(48, 5).Value = 552
$ws.Cells.Item(48, 6).Value = 552000
$ws.Cells.Item(48, 7).Value = 0
$ws.Cells.Item(48, 8).Value = 0
$ws.Cells.Item(49, 5).Value = 613
$ws.Cells.Item(49, 6).Value = 613000
$ws.Cells.Item(49, 7).Value = 0
$ws.Cells.Item(49, 8).Value = 0
$ws.Cells.Item(50, 5).Value = 231
$ws.Cells.Item(50, 6).Value = 231000
$ws.Cells.Item(50, 7).Value = 0
$ws.Cells.Item(50, 8).Value = 0
$ws.Cells.Item(51, 5).Value = 33.75
$ws.Cells.Item(51, 6).Value = 34087.5
$ws.Cells.Item(51, 7).Value = 0
$ws.Cells.Item(51, 8).Value = 0
$ws.Cells.Item(52, 5).Value = 334
$ws.Cells.Item(52, 6).Value = 334000
$ws.Cells.Item(52, 7).Value = 0
$ws.Cells.Item(52, 8).Value = 0
$ws.Cells.Item(53, 5).Value = 250.5
$ws.Cells.Item(53, 6).Value = 250500
$ws.Cells.Item(53, 7).Value = 0
$ws.Cells.Item(53, 8).Value = 0

# Dimension will be recalculated automatically by Excel based on used range
# (A1:H53) once the cells above are populated.
